$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Pagos" (F) and "Inscrições homologadas" (H) columns for the
# affected rows. H = F + G in each case (unchanged G).

$ws.Range("F10").Value = 563
$ws.Range("H10").Value = 658

$ws.Range("F11").Value = 378
$ws.Range("H11").Value = 443

$ws.Range("F12").Value = 612
$ws.Range("H12").Value = 698

$ws.Range("F15").Value = 130
$ws.Range("H15").Value = 181

$ws.Range("F24").Value = 223
$ws.Range("H24").Value = 253

$ws.Range("F25").Value = 254
$ws.Range("H25").Value = 314

$ws.Range("F31").Value = 47
$ws.Range("H31").Value = 74

$ws.Range("F42").Value = 375
$ws.Range("H42").Value = 436
